$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete entire row 9 (NpcMenuFunctionType=4, GoToMyhome), shifting rows 10-23 up
$ws.Rows.Item(9).Delete()

# Update selection to match target state
$ws.Range("C8").Select()
